$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (40 and 41) right after the last existing data row (39),
# copying row 39's formatting (style) down so the new rows match the existing
# sheet's cell style (s="1").
$ws.Rows("39").Copy()
$ws.Rows("40").Insert(-4121)
$ws.Rows("39").Copy()
$ws.Rows("41").Insert(-4121)
$excel.CutCopyMode = 0

# Row 40: Angarsk, female, 2020
$ws.Range("A40").Value = 25703000
$ws.Range("B40").Value = "Ангарский"
$ws.Range("C40").Value = "female"
$ws.Range("D40").Value = 2020
$ws.Range("E40").Value = 0.0585
$ws.Range("F40").Value = 0.06573
$ws.Range("G40").Value = 0.0559
$ws.Range("H40").Value = 0.04904
$ws.Range("I40").Value = 0.04422
$ws.Range("J40").Value = 0.06064
$ws.Range("K40").Value = 0.0869
$ws.Range("L40").Value = 0.0848
$ws.Range("M40").Value = 0.07635
$ws.Range("N40").Value = 0.0729
$ws.Range("O40").Value = 0.12256
$ws.Range("P40").Value = 0.0732
$ws.Range("Q40").Value = 0.0778
$ws.Range("R40").Value = 0.07153

# Row 41: Angarsk, male, 2020
$ws.Range("A41").Value = 25703000
$ws.Range("B41").Value = "Ангарский"
$ws.Range("C41").Value = "male"
$ws.Range("D41").Value = 2020
$ws.Range("E41").Value = 0.0676
$ws.Range("F41").Value = 0.0735
$ws.Range("G41").Value = 0.0637
$ws.Range("H41").Value = 0.05762
$ws.Range("I41").Value = 0.0528
$ws.Range("J41").Value = 0.07214
$ws.Range("K41").Value = 0.10266
$ws.Range("L41").Value = 0.0963
$ws.Range("M41").Value = 0.07806
$ws.Range("N41").Value = 0.0716
$ws.Range("O41").Value = 0.1064
$ws.Range("P41").Value = 0.06085
$ws.Range("Q41").Value = 0.0535
$ws.Range("R41").Value = 0.0431

# Update sheet view selection to mirror the authored change.
$ws.Range("S32").Select()
